# Daily attendance processing - 2025-12-17 17:32:22
#
# Normalizes the "Recorded By" column (G) on the active sheet: whenever the
# first name in the comma-separated list of recorders is "System", move it
# so it no longer leads the list (swap it with the entry right after it).
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#      "System, backup@backdoor.com, system" -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value()

    if ($value -eq $null) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Length -ge 2 -and $parts[0] -eq "System") {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value = [string]::Join(", ", $parts)
    }
}
